# Commit: "Common: Initial form for build"
# Adds 10 new translation rows (r=201..210) to the "Import" sheet, reusing
# the formatting of the last existing data row (r=200), and moves the
# active selection to B196 (as in the target workbook).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Import")

# New rows to append, all in language "cs" (column A), with translation
# key (column B) and Czech text (column C).
$rows = @(
    @("lab.build.title", "Buildy"),
    @("lab.build.subtitle", "V této sekci je možné spravovat a vytvářet buildy."),
    @("lab.build.index.menu", "Obecné"),
    @("lab.build.create.menu", "Nový build"),
    @("lab.build.list.menu", "Seznam buildů"),
    @("lab.build.create.title", "Nový build"),
    @("lab.build.create.subtitle", " "),
    @("lab.build.name.label", "Jméno"),
    @("lab.build.name.label.tooltip", "Jméno buildu musí být unikátní; později bude možné ho použít pro evidenci požitků (vapování)."),
    @("lab.build.atomizerId.label", "Atomizér")
)

$startRow = 201

for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i

    # Clone formatting (style) of the last existing row (200) so the new
    # rows look identical to the existing translation rows.
    $ws.Range("A200:C200").Copy()
    $ws.Range("A$r" + ":C$r").PasteSpecial(-4122)

    $ws.Cells.Item($r, 1).Value = "cs"
    $ws.Cells.Item($r, 2).Value = $rows[$i][0]
    $ws.Cells.Item($r, 3).Value = $rows[$i][1]
}

$excel.CutCopyMode = 0

# Move the active selection, matching the committed workbook state.
$ws.Range("B196").Select() | Out-Null
